# Updated cryptos list on Fri Jul 12 21:35:12 UTC 2024 with GitHub Actions
#
# The Price column (D) is stored as text in the workbook (some values use a
# "." as thousands separator, e.g. "57.587.01", so it can never be real
# numbers). Several updated prices (e.g. "530.42") DO parse as valid
# numbers though, so a plain .Value assignment would silently coerce them
# to a numeric cell. To keep them as text (matching the original layout)
# we momentarily force a text number format, assign the value, then clear
# the format again so no stray style sticks around on the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "57.597.24"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.115.69"
$ws.Range("E3").Value = "  +0.36%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "530.42"
$ws.Range("E5").Value = "  +1.27%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "137.88"
$ws.Range("E6").Value = "  +1.11%  "

# Row 7 - USDC
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Range("D8") "3.113.62"
$ws.Range("E8").Value = "  +0.33%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +5.59%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +0.52%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.64%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.411"
$ws.Range("E12").Value = "  +4.56%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.55%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "3.651.61"
$ws.Range("E14").Value = "  +0.20%  "

# Row 15 - Avalanche
$ws.Range("E15").Value = "  +1.81%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +1.39%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "57.694.04"
$ws.Range("E17").Value = "  +0.37%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "3.120.58"
$ws.Range("E18").Value = "  +0.40%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +2.16%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "12.68"
$ws.Range("E20").Value = "  +2.56%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "360.64"
$ws.Range("E22").Value = "  +4.19%  "

# Row 23 - Dai
Set-TextValue $ws.Range("D23") "0.998"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "69.02"
$ws.Range("E24").Value = "  +2.02%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +0.95%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.08%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.04%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0864"
$ws.Range("E28").Value = "  -3.19%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("E29").Value = "  -1.33%  "

# Row 30 - was PancakeSwap, now RenderToken
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D30") "6.09"
$ws.Range("E30").Value = "  +1.31%  "

# Row 31 - was RenderToken, now PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.87"
$ws.Range("E31").Value = "  +0.14%  "

# Row 32 - EthereumClassic
Set-TextValue $ws.Range("D32") "21.36"
$ws.Range("E32").Value = "  +2.55%  "

# Row 33 - NEARProtocol
Set-TextValue $ws.Range("D33") "5.13"
$ws.Range("E33").Value = "  +4.33%  "

# Row 34 - was Monero, now Fetch.AI
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D34") "1.14"
$ws.Range("E34").Value = "  -0.59%  "

# Row 35 - was Fetch.AI, now Monero
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D35") "159.05"
$ws.Range("E35").Value = "  +0.55%  "

# Row 36 - Aptos
Set-TextValue $ws.Range("D36") "6.06"
$ws.Range("E36").Value = "  +0.20%  "

# Row 37 - ImmutableX
Set-TextValue $ws.Range("D37") "1.29"
$ws.Range("E37").Value = "  +4.76%  "

# Row 38 - EnergySwap
Set-TextValue $ws.Range("D38") "25.51"
$ws.Range("E38").Value = "  -1.21%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +3.78%  "

# Row 40 - Hedera
Set-TextValue $ws.Range("D40") "0.0667"
$ws.Range("E40").Value = "  +1.16%  "

# Row 41 - Maker
Set-TextValue $ws.Range("D41") "2.485.59"
$ws.Range("E41").Value = "  +5.62%  "

# Row 42 - Filecoin
Set-TextValue $ws.Range("D42") "4.02"
$ws.Range("E42").Value = "  -2.98%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  -0.33%  "

# Row 44 - OKB
Set-TextValue $ws.Range("D44") "37.71"
$ws.Range("E44").Value = "  +3.31%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +1.44%  "

# Row 46 - FirstDigitalUSD
Set-TextValue $ws.Range("D46") "0.999"

# Row 47 - ONDO
$ws.Range("E47").Value = "  +2.12%  "

# Row 48 - Cosmos
Set-TextValue $ws.Range("D48") "6.08"
$ws.Range("E48").Value = "  +1.66%  "

# Row 49 - InjectiveProtocol
Set-TextValue $ws.Range("D49") "19.67"
$ws.Range("E49").Value = "  -0.24%  "

# Row 50 - SuiNetwork
Set-TextValue $ws.Range("D50") "0.740"
$ws.Range("E50").Value = "  -2.07%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +2.39%  "
